$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new "Save" header in H1, using the same style as the other header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill H2:H23 with Save flag: 1 if the "sum" column (G) exceeds 10, else 0
for ($row = 2; $row -le 23; $row++) {
    $sumValue = $ws.Cells.Item($row, 7).Value2
    if ($sumValue -gt 10) {
        $ws.Cells.Item($row, 8).Value = 1
    } else {
        $ws.Cells.Item($row, 8).Value = 0
    }
}
